$d = $word.ActiveDocument

$pairs = @(
    @("111×2=", "316×5="),
    @("732×4=", "343×7="),
    @("220×9=", "861×8="),
    @("760×7=", "623×5="),
    @("461×7=", "866×2="),
    @("462×8=", "903×9="),
    @("864×4=", "303×8="),
    @("614×8=", "615×9="),
    @("576×7=", "132×5="),
    @("327×6=", "740×2="),
    @("499×2=", "925×3="),
    @("213×6=", "844×8="),
    @("228×7=", "525×2="),
    @("777×3=", "801×6="),
    @("544×6=", "566×8="),
    @("763×6=", "538×4="),
    @("964×9=", "514×4="),
    @("717×2=", "724×6="),
    @("184×5=", "212×2="),
    @("702×7=", "109×9="),
    @("587×6=", "733×9="),
    @("102×7=", "548×5="),
    @("780×6=", "621×6="),
    @("989×5=", "428×4="),
    @("691×2=", "163×8=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
